$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23.78000000000028
$ws.Range("H2").Value = [double]"1.110223024625157e-16"
$ws.Range("I2").Value = [double]"1.110223024625157e-16"
$ws.Range("L2").Value = 44.00308278221021
$ws.Range("M2").Value = "[37.45874783809205, 50.547417726328376]"
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1.515763422452733
$ws.Range("Q2").Value = "[1.364816027685655, 1.666710817219811]"
$ws.Range("T2").Value = 50.98557058971128
$ws.Range("U2").Value = "[46.743427889954674, 55.227713289467886]"
$ws.Range("X2").Value = 18.04328328328349
$ws.Range("Y2").Value = 17.47199199199219
$ws.Range("Z2").Value = 18.6145745745748
$ws.Range("F3").Value = 23.78000000000028
$ws.Range("L3").Value = 45.01644474860602
$ws.Range("M3").Value = "[38.73763331863638, 51.29525617857566]"
$ws.Range("P3").Value = 1.50318447288881
$ws.Range("Q3").Value = "[1.364816027685655, 1.641552918091965]"
$ws.Range("T3").Value = 49.41764894389998
$ws.Range("U3").Value = "[45.3535333219127, 53.48176456588727]"
$ws.Range("X3").Value = 18.0908908908911
$ws.Range("Y3").Value = 17.56720720720741
$ws.Range("Z3").Value = 18.6145745745748
$ws.Range("F4").Value = 23.78000000000028
$ws.Range("H4").Value = [double]"5.551115123125783e-16"
$ws.Range("I4").Value = [double]"5.551115123125783e-16"
$ws.Range("L4").Value = 49.42167275962659
$ws.Range("M4").Value = "[40.95164488203206, 57.89170063722111]"
$ws.Range("N4").Value = [double]"2.664535259100376e-15"
$ws.Range("O4").Value = [double]"2.664535259100376e-15"
$ws.Range("P4").Value = 1.679289766783733
$ws.Range("Q4").Value = "[1.490605523324886, 1.8679740102425804]"
$ws.Range("T4").Value = 48.78102876590842
$ws.Range("U4").Value = "[43.383409863198885, 54.17864766861796]"
$ws.Range("X4").Value = 17.42438438438459
$ws.Range("Y4").Value = 16.71027027027046
$ws.Range("Z4").Value = 18.13849849849871
$ws.Range("F5").Value = 23.78000000000028
$ws.Range("H5").Value = [double]"1.998401444325282e-15"
$ws.Range("I5").Value = [double]"1.998401444325282e-15"
$ws.Range("L5").Value = 46.81572124644862
$ws.Range("M5").Value = "[37.89731581382784, 55.7341266790694]"
$ws.Range("N5").Value = [double]"8.815170815523743e-14"
$ws.Range("O5").Value = [double]"8.815170815523743e-14"
$ws.Range("P5").Value = 1.805079262422964
$ws.Range("Q5").Value = "[1.591237119836272, 2.0189214050096567]"
$ws.Range("T5").Value = 52.81575688198957
$ws.Range("U5").Value = "[47.53226646115444, 58.0992473028247]"
$ws.Range("X5").Value = 16.9483083083085
$ws.Range("Y5").Value = 16.13897897897917
$ws.Range("Z5").Value = 17.75763763763784
$ws.Range("F6").Value = 23.78000000000028
$ws.Range("H6").Value = [double]"2.109423746787797e-15"
$ws.Range("I6").Value = [double]"2.109423746787797e-15"
$ws.Range("L6").Value = 48.86586261195967
$ws.Range("M6").Value = "[39.60789671454092, 58.12382850937842]"
$ws.Range("N6").Value = [double]"7.394085344003543e-14"
$ws.Range("O6").Value = [double]"7.394085344003543e-14"
$ws.Range("P6").Value = 1.742184514603349
$ws.Range("Q6").Value = "[1.54092132158058, 1.9434477076261185]"
$ws.Range("T6").Value = 53.30803914860124
$ws.Range("U6").Value = "[47.779147460915965, 58.83693083628651]"
$ws.Range("X6").Value = 17.18634634634655
$ws.Range("Y6").Value = 16.42462462462482
$ws.Range("Z6").Value = 17.94806806806828
$ws.Range("F7").Value = 23.78000000000028
$ws.Range("H7").Value = [double]"7.771561172376096e-16"
$ws.Range("I7").Value = [double]"7.771561172376096e-16"
$ws.Range("L7").Value = 44.91161290500328
$ws.Range("M7").Value = "[36.433238906047954, 53.38998690395861]"
$ws.Range("N7").Value = [double]"6.572520305780927e-14"
$ws.Range("O7").Value = [double]"6.572520305780927e-14"
$ws.Range("P7").Value = 1.918289808498272
$ws.Range("Q7").Value = "[1.7170266154755032, 2.1195530015210418]"
$ws.Range("T7").Value = 54.10582091378677
$ws.Range("U7").Value = "[49.169259635147526, 59.04238219242602]"
$ws.Range("X7").Value = 16.51983983984003
$ws.Range("Y7").Value = 15.7581181181183
$ws.Range("Z7").Value = 17.28156156156176
$ws.Range("F8").Value = 23.78000000000028
$ws.Range("H8").Value = [double]"4.440892098500626e-16"
$ws.Range("I8").Value = [double]"4.440892098500626e-16"
$ws.Range("L8").Value = 44.78607887805012
$ws.Range("M8").Value = "[36.34125891820024, 53.2308988379]"
$ws.Range("N8").Value = [double]"6.328271240363392e-14"
$ws.Range("O8").Value = [double]"6.328271240363392e-14"
$ws.Range("P8").Value = 1.86797401024258
$ws.Range("Q8").Value = "[1.666710817219811, 2.0692372032653497]"
$ws.Range("T8").Value = 50.62860508410989
$ws.Range("U8").Value = "[45.78365448051876, 55.47355568770102]"
$ws.Range("X8").Value = 16.71027027027046
$ws.Range("Y8").Value = 15.94854854854873
$ws.Range("Z8").Value = 17.47199199199219
$ws.Range("F9").Value = 23.78000000000028
$ws.Range("H9").Value = [double]"1.942890293094024e-14"
$ws.Range("I9").Value = [double]"1.942890293094024e-14"
$ws.Range("L9").Value = 48.34359851287255
$ws.Range("M9").Value = "[38.53982408776556, 58.14737293797953]"
$ws.Range("N9").Value = [double]"6.441513988875158e-13"
$ws.Range("O9").Value = [double]"6.441513988875158e-13"
$ws.Range("P9").Value = 1.830237161550811
$ws.Range("Q9").Value = "[1.603816069400196, 2.0566582537014266]"
$ws.Range("T9").Value = 55.72356466345199
$ws.Range("U9").Value = "[49.893428807888206, 61.55370051901577]"
$ws.Range("X9").Value = 16.85309309309329
$ws.Range("Y9").Value = 15.99615615615634
$ws.Range("Z9").Value = 17.71003003003023
$ws.Range("F10").Value = 23.78000000000028
$ws.Range("L10").Value = 48.94349420997616
$ws.Range("M10").Value = "[39.26356188419066, 58.62342653576166]"
$ws.Range("N10").Value = [double]"2.930988785010413e-13"
$ws.Range("O10").Value = [double]"2.930988785010413e-13"
$ws.Range("P10").Value = 2.056658253701427
$ws.Range("Q10").Value = "[1.8679740102425804, 2.245342497160273]"
$ws.Range("T10").Value = 55.13066480154065
$ws.Range("U10").Value = "[50.06403508827741, 60.1972945148039]"
$ws.Range("X10").Value = 15.99615615615634
$ws.Range("Y10").Value = 15.28204204204222
$ws.Range("Z10").Value = 16.71027027027046
$ws.Range("F11").Value = 23.71000000000027
$ws.Range("H11").Value = [double]"9.880984919163893e-15"
$ws.Range("I11").Value = [double]"9.880984919163893e-15"
$ws.Range("L11").Value = 47.95236722965571
$ws.Range("M11").Value = "[38.095353211638084, 57.80938124767333]"
$ws.Range("N11").Value = [double]"9.809930645587883e-13"
$ws.Range("O11").Value = [double]"9.809930645587883e-13"
$ws.Range("P11").Value = 1.817658211986887
$ws.Range("Q11").Value = "[1.591237119836272, 2.0440793041375027]"
$ws.Range("T11").Value = 52.41944099093101
$ws.Range("U11").Value = "[46.75303344762549, 58.085848534236526]"
$ws.Range("X11").Value = 16.85095095095114
$ws.Range("Y11").Value = 15.99653653653672
$ws.Range("Z11").Value = 17.70536536536556
$ws.Range("F12").Value = 23.71000000000027
$ws.Range("H12").Value = [double]"1.110223024625157e-16"
$ws.Range("I12").Value = [double]"1.110223024625157e-16"
$ws.Range("L12").Value = 46.18053301041137
$ws.Range("M12").Value = "[37.69265110740975, 54.668414913412995]"
$ws.Range("N12").Value = [double]"2.731148640577885e-14"
$ws.Range("O12").Value = [double]"2.731148640577885e-14"
$ws.Range("P12").Value = 1.918289808498272
$ws.Range("Q12").Value = "[1.7421845146033492, 2.0943951023931957]"
$ws.Range("T12").Value = 47.98543125501081
$ws.Range("U12").Value = "[43.37688984051889, 52.593972669502726]"
$ws.Range("X12").Value = 16.4712112112114
$ws.Range("Y12").Value = 15.80666666666685
$ws.Range("Z12").Value = 17.13575575575595
$ws.Range("F13").Value = 23.71000000000027
$ws.Range("H13").Value = [double]"1.110223024625157e-16"
$ws.Range("I13").Value = [double]"1.110223024625157e-16"
$ws.Range("L13").Value = 48.91035753534821
$ws.Range("M13").Value = "[40.04750399682388, 57.773211073872545]"
$ws.Range("N13").Value = [double]"1.709743457922741e-14"
$ws.Range("O13").Value = [double]"1.709743457922741e-14"
$ws.Range("P13").Value = 1.943447707626119
$ws.Range("Q13").Value = "[1.7547634641672722, 2.132131951084965]"
$ws.Range("T13").Value = 50.6593307614455
$ws.Range("U13").Value = "[45.796098362941386, 55.52256315994962]"
$ws.Range("X13").Value = 16.37627627627646
$ws.Range("Y13").Value = 15.66426426426444
$ws.Range("Z13").Value = 17.08828828828848
$ws.Range("F14").Value = 23.71000000000027
$ws.Range("H14").Value = [double]"3.441691376337985e-15"
$ws.Range("I14").Value = [double]"3.441691376337985e-15"
$ws.Range("L14").Value = 48.77657611107381
$ws.Range("M14").Value = "[38.07295567789517, 59.48019654425244]"
$ws.Range("N14").Value = [double]"7.108980071279802e-12"
$ws.Range("O14").Value = [double]"7.108980071279802e-12"
$ws.Range("P14").Value = 2.106974051957119
$ws.Range("Q14").Value = "[1.8931319093704264, 2.320816194543811]"
$ws.Range("T14").Value = 53.36831310378949
$ws.Range("U14").Value = "[47.76253691752744, 58.97408929005154]"
$ws.Range("X14").Value = 15.75919919919938
$ws.Range("Y14").Value = 14.95225225225242
$ws.Range("Z14").Value = 16.56614614614633
